# Update "想去人数" (want-to-go count) values in columns F on the
# "展览" and "全部类型" worksheets to reflect freshly scraped counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1298
$ws1.Range("F6").Value = 18161
$ws1.Range("F7").Value = 365
$ws1.Range("F10").Value = 6847
$ws1.Range("F18").Value = 1301
$ws1.Range("F19").Value = 228
$ws1.Range("F25").Value = 274
$ws1.Range("F26").Value = 988
$ws1.Range("F33").Value = 12074
$ws1.Range("F38").Value = 3921

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1298
$ws4.Range("F6").Value = 18161
$ws4.Range("F7").Value = 365
$ws4.Range("F10").Value = 6847
$ws4.Range("F18").Value = 1301
$ws4.Range("F19").Value = 228
$ws4.Range("F25").Value = 274
$ws4.Range("F26").Value = 988
$ws4.Range("F35").Value = 12074
$ws4.Range("F40").Value = 3921
